$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 110, shifting existing rows 110-159 down to 111-160.
$ws.Rows.Item(110).Insert()

# Populate the newly inserted row 110 with the new weekly data record.
$ws.Range("A110").Value = 3
$ws.Range("B110").Value = "Femacal de La Calera"
$ws.Range("C110").Value = "Coquimbo"
$ws.Range("D110").Value = 44755
$ws.Range("E110").Value = 5
$ws.Range("F110").Value = 100112026
$ws.Range("G110").Value = "Haba"
$ws.Range("H110").Value = "Sin especificar"
$ws.Range("I110").Value = "Primera"
$ws.Range("J110").Value = 83
$ws.Range("K110").Value = 18000
$ws.Range("L110").Value = 19000
$ws.Range("M110").Value = 18458
$ws.Range("N110").Value = "$/saco 25 kilos"
$ws.Range("O110").Value = "Provincia de Limarí"
$ws.Range("P110").Value = 738
$ws.Range("Q110").Value = 25
$ws.Range("R110").Value = "Hortaliza"
